$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "Muziekgebouw"
$ws.Range("A11").Value = "Centraal Station_A"
$ws.Range("A11").Select()
